# ---------------------------------------------------------------------------
# Re-creates two edits present in the target commit:
#
# 1. The table on slide 6 is switched from the deck's custom "Table_0" style
#    to the built-in "Medium Style 2 - Accent 1" table style
#    ({8B8ADC96-2113-4D87-978A-AF0E2EA414F8}).
#
# 2. The presentation theme's 12-slot colour scheme (ppt/theme/theme1.xml,
#    the theme used by the slide master / all slides) is switched from the
#    "Integral" palette to the default "Office" palette.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------

$tableSlide = $p.Slides.Item(6)
foreach ($shp in $tableSlide.Shapes) {
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{8B8ADC96-2113-4D87-978A-AF0E2EA414F8}", $true)
    }
}

# --- 2. Theme colour scheme -------------------------------------------------
# RGB() packs as 0x00BBGGRR (little-endian BGR) the same way VBA's RGB()
# function / the OLE_COLOR type does.
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office" theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(
    (RGBVal 0x00 0x00 0x00), # dk1
    (RGBVal 0xFF 0xFF 0xFF), # lt1
    (RGBVal 0x44 0x54 0x6A), # dk2
    (RGBVal 0xE7 0xE6 0xE6), # lt2
    (RGBVal 0x5B 0x9B 0xD5), # accent1
    (RGBVal 0xED 0x7D 0x31), # accent2
    (RGBVal 0xA5 0xA5 0xA5), # accent3
    (RGBVal 0xFF 0xC0 0x00), # accent4
    (RGBVal 0x44 0x72 0xC4), # accent5
    (RGBVal 0x70 0xAD 0x47), # accent6
    (RGBVal 0x05 0x63 0xC1), # hlink
    (RGBVal 0x95 0x4F 0x72)  # folHlink
)

$themeSlide = $p.Slides.Item(1)
$tcs = $themeSlide.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
